$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 1 block (rows 5-10) ---
$ws.Range("D6").Value = 1
$ws.Range("H6").Value = 8

$ws.Range("F9").Value = 1
$ws.Range("H9").Value = 10

$ws.Range("H10").Value = 50

# --- Sprint 2 block (rows 11-16) ---
$ws.Range("D12").Value = 1
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 14

$ws.Range("F13").Value = 3
$ws.Range("H13").Value = 14

$ws.Range("H16").Value = 59

# --- Sprint 3 block (rows 17-22) ---
$ws.Range("E18").Value = 2
$ws.Range("H18").Value = 15

$ws.Range("D20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 16

$ws.Range("E21").Value = 3
$ws.Range("H21").Value = 10

$ws.Range("H22").Value = 60

# --- Sprint 4 block (rows 23-28) ---
$ws.Range("E26").Value = 2
$ws.Range("H26").Value = 10

$ws.Range("H28").Value = 51

# --- Grand total row (29) ---
$ws.Range("B29").Value = 35
$ws.Range("C29").Value = 36
$ws.Range("D29").Value = 35
$ws.Range("E29").Value = 35
$ws.Range("F29").Value = 35
$ws.Range("G29").Value = 35

# H29 is a brand-new cell in this row; give it the same look as the other
# total cells on that row (copy format from H28) before writing its value.
$ws.Range("H28").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 211
$excel.CutCopyMode = $false

# --- Sheet view: scroll position and active selection ---
$window = $excel.ActiveWindow
$window.ScrollRow = 14
$window.ScrollColumn = 2
$ws.Range("I32").Select()
